$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.667.64"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.614.71"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "2.629.77"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "3.072.57"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "60.673.62"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "2.623.75"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "2.731.69"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "0.0₃0849"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  +4.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.886"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.91%  "
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.852"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "293.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.624"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +0.24%  "
